{"js": "// Update the worksheet date and every \"axb=\" multiplication prompt to the\n// new values described by the commit diff. Each old string is unique in the\n// document, so an exact, case-sensitive search-and-replace per pair is safe\n// and preserves the existing run formatting (font/size) of each cell.\nconst replacements = [\n  [\"2025-03-06 Thursday\", \"2025-03-07 Friday\"],\n  [\"919\u00d78=\", \"882\u00d75=\"],\n  [\"506\u00d79=\", \"560\u00d79=\"],\n  [\"528\u00d78=\", \"593\u00d79=\"],\n  [\"737\u00d79=\", \"823\u00d76=\"],\n  [\"254\u00d72=\", \"634\u00d72=\"],\n  [\"453\u00d78=\", \"994\u00d79=\"],\n  [\"423\u00d72=\", \"962\u00d78=\"],\n  [\"133\u00d73=\", \"909\u00d74=\"],\n  [\"207\u00d73=\", \"250\u00d76=\"],\n  [\"123\u00d75=\", \"812\u00d79=\"],\n  [\"157\u00d75=\", \"829\u00d77=\"],\n  [\"323\u00d75=\", \"591\u00d74=\"],\n  [\"533\u00d73=\", \"745\u00d77=\"],\n  [\"826\u00d73=\", \"150\u00d78=\"],\n  [\"481\u00d78=\", \"766\u00d79=\"],\n  [\"303\u00d74=\", \"821\u00d74=\"],\n  [\"325\u00d76=\", \"931\u00d72=\"],\n  [\"508\u00d76=\", \"471\u00d76=\"],\n  [\"463\u00d76=\", \"691\u00d79=\"],\n  [\"165\u00d74=\", \"182\u00d78=\"],\n  [\"207\u00d74=\", \"924\u00d75=\"],\n  [\"968\u00d74=\", \"457\u00d73=\"],\n  [\"500\u00d77=\", \"407\u00d73=\"],\n  [\"357\u00d79=\", \"471\u00d78=\"],\n  [\"868\u00d75=\", \"426\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"axb=\" multiplication prompt to the\n# new values described by the commit diff. Each old string is unique in the\n# document, so a Find/Replace (MatchCase, ReplaceAll) per pair is safe and\n# preserves the existing run formatting (font/size) of each cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-03-06 Thursday\", \"2025-03-07 Friday\"),\n    @(\"919\u00d78=\", \"882\u00d75=\"),\n    @(\"506\u00d79=\", \"560\u00d79=\"),\n    @(\"528\u00d78=\", \"593\u00d79=\"),\n    @(\"737\u00d79=\", \"823\u00d76=\"),\n    @(\"254\u00d72=\", \"634\u00d72=\"),\n    @(\"453\u00d78=\", \"994\u00d79=\"),\n    @(\"423\u00d72=\", \"962\u00d78=\"),\n    @(\"133\u00d73=\", \"909\u00d74=\"),\n    @(\"207\u00d73=\", \"250\u00d76=\"),\n    @(\"123\u00d75=\", \"812\u00d79=\"),\n    @(\"157\u00d75=\", \"829\u00d77=\"),\n    @(\"323\u00d75=\", \"591\u00d74=\"),\n    @(\"533\u00d73=\", \"745\u00d77=\"),\n    @(\"826\u00d73=\", \"150\u00d78=\"),\n    @(\"481\u00d78=\", \"766\u00d79=\"),\n    @(\"303\u00d74=\", \"821\u00d74=\"),\n    @(\"325\u00d76=\", \"931\u00d72=\"),\n    @(\"508\u00d76=\", \"471\u00d76=\"),\n    @(\"463\u00d76=\", \"691\u00d79=\"),\n    @(\"165\u00d74=\", \"182\u00d78=\"),\n    @(\"207\u00d74=\", \"924\u00d75=\"),\n    @(\"968\u00d74=\", \"457\u00d73=\"),\n    @(\"500\u00d77=\", \"407\u00d73=\"),\n    @(\"357\u00d79=\", \"471\u00d78=\"),\n    @(\"868\u00d75=\", \"426\u00d73=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $true, $false, $false, $null, $null, $true, $null, $null, $newText, 2)\n}\n"}
